$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove now-unused columns U:AD (previously duplicated data for rows 1-2)
$ws.Range("U1:AD19").Delete() | Out-Null

# Row 1: index header 0..18 in B1:T1 (unchanged values, range just shrank)
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
$ws.Range("G1").Value = 5
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 7
$ws.Range("J1").Value = 8
$ws.Range("K1").Value = 9
$ws.Range("L1").Value = 10
$ws.Range("M1").Value = 11
$ws.Range("N1").Value = 12
$ws.Range("O1").Value = 13
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("R1").Value = 16
$ws.Range("S1").Value = 17
$ws.Range("T1").Value = 18

# Row 2: HKL header row (string labels, reordered for [h,k,l] columns)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "HKL"
$ws.Range("C2").Value = "[2, 1, 1]"
$ws.Range("D2").Value = "[4, 0, 0]"
$ws.Range("E2").Value = "[2, 2, 0]"
$ws.Range("F2").Value = "[2, 0, 0]"
$ws.Range("G2").Value = "[1, 1, 0]"
$ws.Range("H2").Value = "[2, 2, 2]"
$ws.Range("I2").Value = "[3, 1, 0]"
$ws.Range("J2").Value = "[3, 2, 1]"
$ws.Range("K2").Value = "1Pair-A"
$ws.Range("L2").Value = "1Pair-B"
$ws.Range("M2").Value = "2Pairs-A"
$ws.Range("N2").Value = "2Pairs-B"
$ws.Range("O2").Value = "3Pairs-A"
$ws.Range("P2").Value = "3Pairs-B"
$ws.Range("Q2").Value = "3Pairs-C"
$ws.Range("R2").Value = "4Pairs"
$ws.Range("S2").Value = "5A4F"
$ws.Range("T2").Value = "MaxUnique"

# Data rows 3-23 (rows 16-19 become the new "Holden" series; rows 20-23 are the
# HexGrid series shifted down from their old row numbers 16-19)
# Row 3: BT8Hex_2.5
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "BT8Hex_2.5"
$ws.Range("C3").Value = 0.9951878957092909
$ws.Range("D3").Value = 0.9993041966129168
$ws.Range("E3").Value = 1.025144789947068
$ws.Range("F3").Value = 0.9993041966129168
$ws.Range("G3").Value = 1.025144789947068
$ws.Range("H3").Value = 0.9836838572145354
$ws.Range("I3").Value = 1.000447423893858
$ws.Range("J3").Value = 0.9973346660181278
$ws.Range("K3").Value = 1.025144789947068
$ws.Range("L3").Value = 0.9951878957092909
$ws.Range("M3").Value = 0.9972460461611039
$ws.Range("N3").Value = 0.9972460461611039
$ws.Range("O3").Value = 0.9983131720720219
$ws.Range("P3").Value = 1.006545627423092
$ws.Range("Q3").Value = 1.006545627423092
$ws.Range("R3").Value = 1.011195418054086
$ws.Range("S3").Value = 1.011195418054086
$ws.Range("T3").Value = 1.000183804899299

# Row 4: BT8Hex_5
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "BT8Hex_5"
$ws.Range("C4").Value = 0.9896331371669096
$ws.Range("D4").Value = 0.9976974056591632
$ws.Range("E4").Value = 1.049340798865434
$ws.Range("F4").Value = 0.9976974056591632
$ws.Range("G4").Value = 1.049340798865434
$ws.Range("H4").Value = 0.9684263341830096
$ws.Range("I4").Value = 1.00073550742747
$ws.Range("J4").Value = 0.9952914626482522
$ws.Range("K4").Value = 1.049340798865434
$ws.Range("L4").Value = 0.9896331371669096
$ws.Range("M4").Value = 0.9936652714130364
$ws.Range("N4").Value = 0.9936652714130364
$ws.Range("O4").Value = 0.9960220167511812
$ws.Range("P4").Value = 1.012223780563836
$ws.Range("Q4").Value = 1.012223780563835
$ws.Range("R4").Value = 1.021503035139235
$ws.Range("S4").Value = 1.021503035139235
$ws.Range("T4").Value = 1.000187440991706

# Row 5: BT8Hex_10
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "BT8Hex_10"
$ws.Range("C5").Value = 0.9819855445115087
$ws.Range("D5").Value = 0.9997422610411486
$ws.Range("E5").Value = 1.095987374224819
$ws.Range("F5").Value = 0.9997422610411486
$ws.Range("G5").Value = 1.095987374224819
$ws.Range("H5").Value = 0.9446771174439034
$ws.Range("I5").Value = 1.002716086493868
$ws.Range("J5").Value = 0.9922935145957116
$ws.Range("K5").Value = 1.095987374224819
$ws.Range("L5").Value = 0.9819855445115087
$ws.Range("M5").Value = 0.9908639027763286
$ws.Range("N5").Value = 0.9908639027763286
$ws.Range("O5").Value = 0.9948146306821749
$ws.Range("P5").Value = 1.025905059925825
$ws.Range("Q5").Value = 1.025905059925825
$ws.Range("R5").Value = 1.043425638500574
$ws.Range("S5").Value = 1.043425638500574
$ws.Range("T5").Value = 1.00290031638516

# Row 6: BT8Hex_15
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "BT8Hex_15"
$ws.Range("C6").Value = 0.9994216394135198
$ws.Range("D6").Value = 1.015167821135583
$ws.Range("E6").Value = 1.156346822551598
$ws.Range("F6").Value = 1.015167821135583
$ws.Range("G6").Value = 1.156346822551598
$ws.Range("H6").Value = 0.9205057964863077
$ws.Range("I6").Value = 1.009150855697621
$ws.Range("J6").Value = 0.9768733072209441
$ws.Range("K6").Value = 1.156346822551598
$ws.Range("L6").Value = 0.9994216394135198
$ws.Range("M6").Value = 1.007294730274552
$ws.Range("N6").Value = 1.007294730274552
$ws.Range("O6").Value = 1.007913438748908
$ws.Range("P6").Value = 1.056978761033567
$ws.Range("Q6").Value = 1.056978761033567
$ws.Range("R6").Value = 1.081820776413075
$ws.Range("S6").Value = 1.081820776413075
$ws.Range("T6").Value = 1.012911040417596

# Row 7: Spiral2.5
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Spiral2.5"
$ws.Range("C7").Value = 1.00041337334096
$ws.Range("D7").Value = 1.001964744058554
$ws.Range("E7").Value = 0.9999096111591124
$ws.Range("F7").Value = 1.001964744058554
$ws.Range("G7").Value = 0.9999096111591124
$ws.Range("H7").Value = 1.001764394695247
$ws.Range("I7").Value = 0.9996428419171198
$ws.Range("J7").Value = 0.9996222708504067
$ws.Range("K7").Value = 0.9999096111591124
$ws.Range("L7").Value = 1.00041337334096
$ws.Range("M7").Value = 1.001189058699757
$ws.Range("N7").Value = 1.001189058699757
$ws.Range("O7").Value = 1.000673653105544
$ws.Range("P7").Value = 1.000762576186209
$ws.Range("Q7").Value = 1.000762576186209
$ws.Range("R7").Value = 1.000549334929435
$ws.Range("S7").Value = 1.000549334929435
$ws.Range("T7").Value = 1.000552872670233

# Row 8: Spiral5
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Spiral5"
$ws.Range("C8").Value = 1.000736091982743
$ws.Range("D8").Value = 1.003752402361496
$ws.Range("E8").Value = 1.00159975389846
$ws.Range("F8").Value = 1.003752402361496
$ws.Range("G8").Value = 1.00159975389846
$ws.Range("H8").Value = 1.003999889073205
$ws.Range("I8").Value = 0.999033185473417
$ws.Range("J8").Value = 0.9990214260114924
$ws.Range("K8").Value = 1.00159975389846
$ws.Range("L8").Value = 1.000736091982743
$ws.Range("M8").Value = 1.00224424717212
$ws.Range("N8").Value = 1.00224424717212
$ws.Range("O8").Value = 1.001173893272552
$ws.Range("P8").Value = 1.0020294160809
$ws.Range("Q8").Value = 1.0020294160809
$ws.Range("R8").Value = 1.00192200053529
$ws.Range("S8").Value = 1.00192200053529
$ws.Range("T8").Value = 1.001357124800136

# Row 9: Spiral7.5
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Spiral7.5"
$ws.Range("C9").Value = 1.000241581506902
$ws.Range("D9").Value = 1.004443164042033
$ws.Range("E9").Value = 1.005422243409984
$ws.Range("F9").Value = 1.004443164042033
$ws.Range("G9").Value = 1.005422243409984
$ws.Range("H9").Value = 1.003644768235944
$ws.Range("I9").Value = 0.9985767206637098
$ws.Range("J9").Value = 0.998626064786322
$ws.Range("K9").Value = 1.005422243409984
$ws.Range("L9").Value = 1.000241581506902
$ws.Range("M9").Value = 1.002342372774467
$ws.Range("N9").Value = 1.002342372774467
$ws.Range("O9").Value = 1.001087155404215
$ws.Range("P9").Value = 1.003368996319639
$ws.Range("Q9").Value = 1.003368996319639
$ws.Range("R9").Value = 1.003882308092226
$ws.Range("S9").Value = 1.003882308092226
$ws.Range("T9").Value = 1.001825757107482

# Row 10: Spiral10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Spiral10"
$ws.Range("C10").Value = 1.0014551184254
$ws.Range("D10").Value = 1.010662248990591
$ws.Range("E10").Value = 1.007755921168244
$ws.Range("F10").Value = 1.010662248990591
$ws.Range("G10").Value = 1.007755921168244
$ws.Range("H10").Value = 1.010666914911563
$ws.Range("I10").Value = 0.9969642054895369
$ws.Range("J10").Value = 0.9970236676244635
$ws.Range("K10").Value = 1.007755921168244
$ws.Range("L10").Value = 1.0014551184254
$ws.Range("M10").Value = 1.006058683707996
$ws.Range("N10").Value = 1.006058683707996
$ws.Range("O10").Value = 1.00302719096851
$ws.Range("P10").Value = 1.006624429528079
$ws.Range("Q10").Value = 1.006624429528078
$ws.Range("R10").Value = 1.00690730243812
$ws.Range("S10").Value = 1.00690730243812
$ws.Range("T10").Value = 1.0040880127683

# Row 11: Spiral15
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral15"
$ws.Range("C11").Value = 0.9994498954210648
$ws.Range("D11").Value = 1.013413964427211
$ws.Range("E11").Value = 1.023786138663091
$ws.Range("F11").Value = 1.013413964427211
$ws.Range("G11").Value = 1.023786138663091
$ws.Range("H11").Value = 1.012179255866814
$ws.Range("I11").Value = 0.9955104353890499
$ws.Range("J11").Value = 0.9949875146238187
$ws.Range("K11").Value = 1.023786138663091
$ws.Range("L11").Value = 0.9994498954210648
$ws.Range("M11").Value = 1.006431929924138
$ws.Range("N11").Value = 1.006431929924138
$ws.Range("O11").Value = 1.002791431745775
$ws.Range("P11").Value = 1.012216666170455
$ws.Range("Q11").Value = 1.012216666170455
$ws.Range("R11").Value = 1.015109034293614
$ws.Range("S11").Value = 1.015109034293614
$ws.Range("T11").Value = 1.006554534065175

# Row 12: OffsetF45
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "OffsetF45"
$ws.Range("C12").Value = 0.9058000676625266
$ws.Range("D12").Value = 0.5242975229237483
$ws.Range("E12").Value = 1.367451572519013
$ws.Range("F12").Value = 0.5242975229237483
$ws.Range("G12").Value = 1.367451572519013
$ws.Range("H12").Value = 0.2955190547006418
$ws.Range("I12").Value = 1.10415844838367
$ws.Range("J12").Value = 1.047797231911656
$ws.Range("K12").Value = 1.367451572519013
$ws.Range("L12").Value = 0.9058000676625266
$ws.Range("M12").Value = 0.7150487952931375
$ws.Range("N12").Value = 0.7150487952931375
$ws.Range("O12").Value = 0.8447520129899818
$ws.Range("P12").Value = 0.9325163877017628
$ws.Range("Q12").Value = 0.9325163877017628
$ws.Range("R12").Value = 1.041250183906075
$ws.Range("S12").Value = 1.041250183906075
$ws.Range("T12").Value = 0.8741706496835429

# Row 13: OffsetA45
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "OffsetA45"
$ws.Range("C13").Value = 0.9559112217430918
$ws.Range("D13").Value = 1.465386549296405
$ws.Range("E13").Value = 0.9135524736007585
$ws.Range("F13").Value = 1.465386549296405
$ws.Range("G13").Value = 0.9135524736007585
$ws.Range("H13").Value = 1.291586212764821
$ws.Range("I13").Value = 0.9487592721989307
$ws.Range("J13").Value = 0.9498782134979767
$ws.Range("K13").Value = 0.9135524736007585
$ws.Range("L13").Value = 0.9559112217430918
$ws.Range("M13").Value = 1.210648885519749
$ws.Range("N13").Value = 1.210648885519749
$ws.Range("O13").Value = 1.123352347746143
$ws.Range("P13").Value = 1.111616748213419
$ws.Range("Q13").Value = 1.111616748213419
$ws.Range("R13").Value = 1.062100679560254
$ws.Range("S13").Value = 1.062100679560254
$ws.Range("T13").Value = 1.087512323850331

# Row 14: OffsetFTD
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "OffsetFTD"
$ws.Range("C14").Value = 0.9742458777296152
$ws.Range("D14").Value = 1.737560417288984
$ws.Range("E14").Value = 0.6255072321765257
$ws.Range("F14").Value = 1.737560417288984
$ws.Range("G14").Value = 0.6255072321765257
$ws.Range("H14").Value = 1.314279643418631
$ws.Range("I14").Value = 1.148330836567969
$ws.Range("J14").Value = 0.9081754660064649
$ws.Range("K14").Value = 0.6255072321765257
$ws.Range("L14").Value = 0.9742458777296152
$ws.Range("M14").Value = 1.3559031475093
$ws.Range("N14").Value = 1.3559031475093
$ws.Range("O14").Value = 1.286712377195523
$ws.Range("P14").Value = 1.112437842398375
$ws.Range("Q14").Value = 1.112437842398375
$ws.Range("R14").Value = 0.9907051898429127
$ws.Range("S14").Value = 0.9907051898429127
$ws.Range("T14").Value = 1.118016578864698

# Row 15: OffsetATD
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "OffsetATD"
$ws.Range("C15").Value = 1.016310253986152
$ws.Range("D15").Value = 0.6323739763359562
$ws.Range("E15").Value = 0.97840586922342
$ws.Range("F15").Value = 0.6323739763359562
$ws.Range("G15").Value = 0.97840586922342
$ws.Range("H15").Value = 0.7473766971540057
$ws.Range("I15").Value = 1.070986411344868
$ws.Range("J15").Value = 1.00175823847852
$ws.Range("K15").Value = 0.97840586922342
$ws.Range("L15").Value = 1.016310253986152
$ws.Range("M15").Value = 0.8243421151610543
$ws.Range("N15").Value = 0.8243421151610543
$ws.Range("O15").Value = 0.9065568805556588
$ws.Range("P15").Value = 0.8756966998485095
$ws.Range("Q15").Value = 0.8756966998485095
$ws.Range("R15").Value = 0.9013739921922371
$ws.Range("S15").Value = 0.9013739921922371
$ws.Range("T15").Value = 0.9078685744204869

# Row 16: Holden2.5
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("C16").Value = 0.9090269504671833
$ws.Range("D16").Value = 1.003260044362084
$ws.Range("E16").Value = 1.412039083785934
$ws.Range("F16").Value = 1.003260044362084
$ws.Range("G16").Value = 1.412039083785934
$ws.Range("H16").Value = 0.7055446589074096
$ws.Range("I16").Value = 1.017569624230751
$ws.Range("J16").Value = 0.9593735847823748
$ws.Range("K16").Value = 1.412039083785934
$ws.Range("L16").Value = 0.9090269504671833
$ws.Range("M16").Value = 0.9561434974146337
$ws.Range("N16").Value = 0.9561434974146337
$ws.Range("O16").Value = 0.9766188730200062
$ws.Range("P16").Value = 1.108108692871734
$ws.Range("Q16").Value = 1.108108692871734
$ws.Range("R16").Value = 1.184091290600284
$ws.Range("S16").Value = 1.184091290600284
$ws.Range("T16").Value = 1.001135657755956

# Row 17: Holden5
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "Holden5"
$ws.Range("C17").Value = 0.9528361611472602
$ws.Range("D17").Value = 1.061178117595809
$ws.Range("E17").Value = 1.237393226341469
$ws.Range("F17").Value = 1.061178117595809
$ws.Range("G17").Value = 1.237393226341469
$ws.Range("H17").Value = 0.852760247639252
$ws.Range("I17").Value = 1.013167669915372
$ws.Range("J17").Value = 0.9636308759124764
$ws.Range("K17").Value = 1.237393226341469
$ws.Range("L17").Value = 0.9528361611472602
$ws.Range("M17").Value = 1.007007139371535
$ws.Range("N17").Value = 1.007007139371535
$ws.Range("O17").Value = 1.009060649552813
$ws.Range("P17").Value = 1.083802501694846
$ws.Range("Q17").Value = 1.083802501694846
$ws.Range("R17").Value = 1.122200182856501
$ws.Range("S17").Value = 1.122200182856501
$ws.Range("T17").Value = 1.013494383091939

# Row 18: Holden10
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "Holden10"
$ws.Range("C18").Value = 1.038881895556477
$ws.Range("D18").Value = 1.173449173555999
$ws.Range("E18").Value = 0.8943679930788919
$ws.Range("F18").Value = 1.173449173555999
$ws.Range("G18").Value = 0.8943679930788919
$ws.Range("H18").Value = 1.143282015938244
$ws.Range("I18").Value = 1.004136803828313
$ws.Range("J18").Value = 0.9722756318678779
$ws.Range("K18").Value = 0.8943679930788919
$ws.Range("L18").Value = 1.038881895556477
$ws.Range("M18").Value = 1.106165534556238
$ws.Range("N18").Value = 1.106165534556238
$ws.Range("O18").Value = 1.07215595764693
$ws.Range("P18").Value = 1.035566354063789
$ws.Range("Q18").Value = 1.035566354063789
$ws.Range("R18").Value = 1.000266763817565
$ws.Range("S18").Value = 1.000266763817565
$ws.Range("T18").Value = 1.0377322523043

# Row 19: Holden15
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "Holden15"
$ws.Range("C19").Value = 1.060055515148158
$ws.Range("D19").Value = 1.238894755146833
$ws.Range("E19").Value = 0.8409583458274509
$ws.Range("F19").Value = 1.238894755146833
$ws.Range("G19").Value = 0.8409583458274509
$ws.Range("H19").Value = 1.224218702463185
$ws.Range("I19").Value = 0.9957690713245142
$ws.Range("J19").Value = 0.9649188822764647
$ws.Range("K19").Value = 0.8409583458274509
$ws.Range("L19").Value = 1.060055515148158
$ws.Range("M19").Value = 1.149475135147495
$ws.Range("N19").Value = 1.149475135147495
$ws.Range("O19").Value = 1.098239780539835
$ws.Range("P19").Value = 1.046636205374147
$ws.Range("Q19").Value = 1.046636205374147
$ws.Range("R19").Value = 0.995216740487473
$ws.Range("S19").Value = 0.995216740487473
$ws.Range("T19").Value = 1.054135878697767

# Row 20: HexGrid-90degTilt2.5degRes
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "HexGrid-90degTilt2.5degRes"
$ws.Range("C20").Value = 1.003054716196678
$ws.Range("D20").Value = 1.000547447310705
$ws.Range("E20").Value = 0.993455293527991
$ws.Range("F20").Value = 1.000547447310705
$ws.Range("G20").Value = 0.993455293527991
$ws.Range("H20").Value = 1.00521871669402
$ws.Range("I20").Value = 1.000014457200248
$ws.Range("J20").Value = 0.9998963774510461
$ws.Range("K20").Value = 0.993455293527991
$ws.Range("L20").Value = 1.003054716196678
$ws.Range("M20").Value = 1.001801081753692
$ws.Range("N20").Value = 1.001801081753692
$ws.Range("O20").Value = 1.001205540235877
$ws.Range("P20").Value = 0.999019152345125
$ws.Range("Q20").Value = 0.9990191523451247
$ws.Range("R20").Value = 0.9976281876408413
$ws.Range("S20").Value = 0.9976281876408413
$ws.Range("T20").Value = 1.000364501396781

# Row 21: HexGrid-90degTilt5degRes
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C21").Value = 0.9963540467058953
$ws.Range("D21").Value = 1.001563704069303
$ws.Range("E21").Value = 1.008366134383959
$ws.Range("F21").Value = 1.001563704069303
$ws.Range("G21").Value = 1.008366134383959
$ws.Range("H21").Value = 0.9957487340508335
$ws.Range("I21").Value = 0.9994552097983425
$ws.Range("J21").Value = 0.9995709910448528
$ws.Range("K21").Value = 1.008366134383959
$ws.Range("L21").Value = 0.9963540467058953
$ws.Range("M21").Value = 0.9989588753875992
$ws.Range("N21").Value = 0.9989588753875992
$ws.Range("O21").Value = 0.9991243201911804
$ws.Range("P21").Value = 1.002094628386386
$ws.Range("Q21").Value = 1.002094628386386
$ws.Range("R21").Value = 1.003662504885779
$ws.Range("S21").Value = 1.003662504885779
$ws.Range("T21").Value = 1.000176470008864

# Row 22: HexGrid-90degTilt10degRes
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "HexGrid-90degTilt10degRes"
$ws.Range("C22").Value = 1.004935907958623
$ws.Range("D22").Value = 0.9938723881077177
$ws.Range("E22").Value = 1.003791802394072
$ws.Range("F22").Value = 0.9938723881077177
$ws.Range("G22").Value = 1.003791802394072
$ws.Range("H22").Value = 1.009178850053638
$ws.Range("I22").Value = 0.9994038467739389
$ws.Range("J22").Value = 0.9977174400839748
$ws.Range("K22").Value = 1.003791802394072
$ws.Range("L22").Value = 1.004935907958623
$ws.Range("M22").Value = 0.9994041480331703
$ws.Range("N22").Value = 0.9994041480331703
$ws.Range("O22").Value = 0.9994040476134266
$ws.Range("P22").Value = 1.000866699486805
$ws.Range("Q22").Value = 1.000866699486805
$ws.Range("R22").Value = 1.001597975213621
$ws.Range("S22").Value = 1.001597975213621
$ws.Range("T22").Value = 1.001483372561994

# Row 23: HexGrid-90degTilt15degRes
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "HexGrid-90degTilt15degRes"
$ws.Range("C23").Value = 1.022626762243706
$ws.Range("D23").Value = 0.9971050077671598
$ws.Range("E23").Value = 0.9495585098975111
$ws.Range("F23").Value = 0.9971050077671598
$ws.Range("G23").Value = 0.9495585098975111
$ws.Range("H23").Value = 1.023269698706081
$ws.Range("I23").Value = 1.001910314625178
$ws.Range("J23").Value = 1.004488788656324
$ws.Range("K23").Value = 0.9495585098975111
$ws.Range("L23").Value = 1.022626762243706
$ws.Range("M23").Value = 1.009865885005433
$ws.Range("N23").Value = 1.009865885005433
$ws.Range("O23").Value = 1.007214028212015
$ws.Range("P23").Value = 0.9897634266361256
$ws.Range("Q23").Value = 0.9897634266361256
$ws.Range("R23").Value = 0.979712197451472
$ws.Range("S23").Value = 0.979712197451472
$ws.Range("T23").Value = 0.9998265136493267

# New rows 20-23 need the same bold/centered/bordered style that column A already
# carries for rows 2-19; copy it over (format-only paste) rather than rebuilding it
# attribute-by-attribute so no redundant style entries get created.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A20:A23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false